# EPICP-1: changed unit from mg/d to g/d for sodium and potassium intake
# in DD_EPICP_INES (Variables sheet).
#
# Row 25 (variable "mna"): label changes from
#   "natrium intake at baseline [mg/d]" -> "sodium intake at baseline [g/d]"
# Row 26 (variable "mk"): label changes from
#   "potassium intake at baseline [mg/d]" -> "potassium intake at baseline [g/d]"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

$ws.Range("C25").Value = "sodium intake at baseline [g/d]"
$ws.Range("C26").Value = "potassium intake at baseline [g/d]"

# Update the active selection to match the edited cell.
$ws.Activate()
$ws.Range("C25").Select()
